$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 319.7143
$ws.Range("J33").Value = 130
$ws.Range("L33").Value = 130
$ws.Range("N33").Value = -588
$ws.Range("H43").Value = 3105.7693
$ws.Range("J43").Value = 2208
$ws.Range("L43").Value = 2208
$ws.Range("N43").Value = -2346
$ws.Range("H53").Value = 1384.4286
$ws.Range("J53").Value = 1368
$ws.Range("L53").Value = 1368
$ws.Range("N53").Value = -2642
$ws.Range("H64").Value = 4709.143
$ws.Range("I64").Value = 4094.7778
$ws.Range("K64").Value = 4094.7778
$ws.Range("M64").Value = -3846.7778
$ws.Range("H67").Value = 4709.143
$ws.Range("I67").Value = 4094.7778
$ws.Range("K67").Value = 4094.7778
$ws.Range("M67").Value = -3236.7778
$ws.Range("H82").Value = 25950
$ws.Range("I82").Value = 22678.666
$ws.Range("K82").Value = 68035.99800000001
$ws.Range("M82").Value = -67629.99800000001
$ws.Range("H85").Value = 25950
$ws.Range("I85").Value = 22678.666
$ws.Range("K85").Value = 68035.99800000001
$ws.Range("M85").Value = -66631.99800000001
$ws.Range("H106").Value = 3451.8333
$ws.Range("I106").Value = 3744.2
$ws.Range("K106").Value = 3744.2
$ws.Range("M106").Value = -3113.2
$ws.Range("H138").Value = 3712.9834
$ws.Range("I138").Value = 2830.5789
$ws.Range("J138").Value = 4121.9023
$ws.Range("K138").Value = 8491.736699999999
$ws.Range("L138").Value = 12365.7069
$ws.Range("M138").Value = -3351.736699999999
$ws.Range("N138").Value = -22645.7069
$ws.Range("H141").Value = 20838926
$ws.Range("I141").Value = 22730382
$ws.Range("K141").Value = 68191146
$ws.Range("M141").Value = -68185966

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5489.4653
$ws.Range("I32").Value = 5506.9644
$ws.Range("J32").Value = 4999.5
$ws.Range("K32").Value = 5506.9644
$ws.Range("L32").Value = 4999.5
$ws.Range("M32").Value = -5219.9644
$ws.Range("N32").Value = -5573.5
$ws.Range("H45").Value = 4966.8823
$ws.Range("I45").Value = 2039.1428
$ws.Range("J45").Value = 7016.3
$ws.Range("K45").Value = 2039.1428
$ws.Range("L45").Value = 7016.3
$ws.Range("M45").Value = -1662.1428
$ws.Range("N45").Value = -7770.3
$ws.Range("H60").Value = 70591.57000000001
$ws.Range("I60").Value = 70591.57000000001
$ws.Range("K60").Value = 70591.57000000001
$ws.Range("M60").Value = -69858.57000000001
$ws.Range("H61").Value = 12913416
$ws.Range("I61").Value = 15221456
$ws.Range("J61").Value = 1116766.5
$ws.Range("K61").Value = 15221456
$ws.Range("L61").Value = 1116766.5
$ws.Range("M61").Value = -15221244
$ws.Range("N61").Value = -1117190.5
$ws.Range("H110").Value = 5642.7144
$ws.Range("I110").Value = 6068.263
$ws.Range("K110").Value = 6068.263
$ws.Range("M110").Value = -4023.263
$ws.Range("H124").Value = 18721.166
$ws.Range("J124").Value = 18721.166
$ws.Range("L124").Value = 18721.166
$ws.Range("N124").Value = -28541.166
$ws.Range("H132").Value = 2003463.2
$ws.Range("I132").Value = 3303.186
$ws.Range("J132").Value = 14290161
$ws.Range("K132").Value = 9909.558000000001
$ws.Range("L132").Value = 42870483
$ws.Range("M132").Value = -7379.558000000001
$ws.Range("N132").Value = -42875543
$ws.Range("H136").Value = 12913416
$ws.Range("I136").Value = 15221456
$ws.Range("J136").Value = 1116766.5
$ws.Range("K136").Value = 45664368
$ws.Range("L136").Value = 3350299.5
$ws.Range("M136").Value = -45661818
$ws.Range("N136").Value = -3355399.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 917303.9399999999
$ws.Range("I105").Value = 1633513.1
$ws.Range("J105").Value = 5765
$ws.Range("K105").Value = 1633513.1
$ws.Range("L105").Value = 5765
$ws.Range("M105").Value = -1631766.1
$ws.Range("N105").Value = -9259
$ws.Range("H134").Value = 2274438.8
$ws.Range("I134").Value = 1672.2572
$ws.Range("J134").Value = 11112975
$ws.Range("K134").Value = 5016.7716
$ws.Range("L134").Value = 33338925
$ws.Range("M134").Value = -2481.7716
$ws.Range("N134").Value = -33343995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26318848
$ws.Range("I31").Value = 41669356
$ws.Range("J31").Value = 3692.7144
$ws.Range("K31").Value = 41669356
$ws.Range("L31").Value = 3692.7144
$ws.Range("M31").Value = -41669061
$ws.Range("N31").Value = -4282.7144
$ws.Range("H34").Value = 26318848
$ws.Range("I34").Value = 41669356
$ws.Range("J34").Value = 3692.7144
$ws.Range("K34").Value = 41669356
$ws.Range("L34").Value = 3692.7144
$ws.Range("M34").Value = -41669154
$ws.Range("N34").Value = -4096.7144
$ws.Range("H36").Value = 60016
$ws.Range("I36").Value = 55024
$ws.Range("K36").Value = 55024
$ws.Range("M36").Value = -54636
$ws.Range("H40").Value = 60016
$ws.Range("I40").Value = 55024
$ws.Range("K40").Value = 55024
$ws.Range("M40").Value = -54864
$ws.Range("H58").Value = 3001.3635
$ws.Range("I58").Value = 2847.5386
$ws.Range("J58").Value = 3223.5557
$ws.Range("K58").Value = 2847.5386
$ws.Range("L58").Value = 3223.5557
$ws.Range("M58").Value = -2644.5386
$ws.Range("N58").Value = -3629.5557
$ws.Range("H86").Value = 12658.833
$ws.Range("I86").Value = 15235.5
$ws.Range("K86").Value = 15235.5
$ws.Range("M86").Value = -14112.5
$ws.Range("H89").Value = 12658.833
$ws.Range("I89").Value = 15235.5
$ws.Range("K89").Value = 76177.5
$ws.Range("M89").Value = -70561.5
$ws.Range("H99").Value = 9963.697
$ws.Range("I99").Value = 7693.125
$ws.Range("J99").Value = 12100.706
$ws.Range("K99").Value = 7693.125
$ws.Range("L99").Value = 12100.706
$ws.Range("M99").Value = -6195.125
$ws.Range("N99").Value = -15096.706
$ws.Range("H122").Value = 2716.25
$ws.Range("I122").Value = 2649.5334
$ws.Range("K122").Value = 7948.600199999999
$ws.Range("M122").Value = -5498.600199999999
$ws.Range("H126").Value = 9963.697
$ws.Range("I126").Value = 7693.125
$ws.Range("J126").Value = 12100.706
$ws.Range("K126").Value = 23079.375
$ws.Range("L126").Value = 36302.118
$ws.Range("M126").Value = -20609.375
$ws.Range("N126").Value = -41242.118
$ws.Range("H136").Value = 3001.3635
$ws.Range("I136").Value = 2847.5386
$ws.Range("J136").Value = 3223.5557
$ws.Range("K136").Value = 8542.6158
$ws.Range("L136").Value = 9670.667099999999
$ws.Range("M136").Value = -5992.6158
$ws.Range("N136").Value = -14770.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 4205.4614
$ws.Range("I12").Value = 107.75
$ws.Range("K12").Value = 323.25
$ws.Range("M12").Value = -150.25
$ws.Range("H131").Value = 4531.3184
$ws.Range("I131").Value = 2972.4167
$ws.Range("K131").Value = 8917.250100000001
$ws.Range("M131").Value = -3877.250100000001
$ws.Range("H137").Value = 11080.714
$ws.Range("I137").Value = 6400
$ws.Range("J137").Value = 14591.25
$ws.Range("K137").Value = 19200
$ws.Range("L137").Value = 43773.75
$ws.Range("M137").Value = -14100
$ws.Range("N137").Value = -53973.75
$ws.Range("H138").Value = 12715.538
$ws.Range("J138").Value = 14443.667
$ws.Range("L138").Value = 43331.001
$ws.Range("N138").Value = -53611.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 99997
$ws.Range("J130").Value = 99997
$ws.Range("L130").Value = 99997
$ws.Range("N130").Value = -110037
$ws.Range("H131").Value = 86972
$ws.Range("J131").Value = 86972
$ws.Range("L131").Value = 86972
$ws.Range("N131").Value = -97052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1585.5555
$ws.Range("I16").Value = 638.0606
$ws.Range("J16").Value = 12008
$ws.Range("K16").Value = 638.0606
$ws.Range("L16").Value = 12008
$ws.Range("M16").Value = -468.0606
$ws.Range("N16").Value = -12348
$ws.Range("H22").Value = 767.8570999999999
$ws.Range("I22").Value = 990.25
$ws.Range("K22").Value = 990.25
$ws.Range("M22").Value = -695.25
$ws.Range("H27").Value = 767.8570999999999
$ws.Range("I27").Value = 990.25
$ws.Range("K27").Value = 990.25
$ws.Range("M27").Value = -883.25
$ws.Range("H40").Value = 5186.9565
$ws.Range("I40").Value = 5161.7617
$ws.Range("K40").Value = 5161.7617
$ws.Range("M40").Value = -5025.7617
$ws.Range("H132").Value = 3483.4783
$ws.Range("I132").Value = 1884
$ws.Range("J132").Value = 6482.5
$ws.Range("K132").Value = 5652
$ws.Range("L132").Value = 19447.5
$ws.Range("M132").Value = -3122
$ws.Range("N132").Value = -24507.5
$ws.Range("H136").Value = 3764.5
$ws.Range("I136").Value = 3611.2917
$ws.Range("J136").Value = 4683.75
$ws.Range("K136").Value = 10833.8751
$ws.Range("L136").Value = 14051.25
$ws.Range("M136").Value = -8283.875100000001
$ws.Range("N136").Value = -19151.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2883.9
$ws.Range("J96").Value = 2846.3333
$ws.Range("L96").Value = 2846.3333
$ws.Range("N96").Value = -5592.3333
$ws.Range("H126").Value = 2408.853
$ws.Range("I126").Value = 2663.8635
$ws.Range("J126").Value = 1941.3334
$ws.Range("K126").Value = 7991.5905
$ws.Range("L126").Value = 5824.0002
$ws.Range("M126").Value = -5521.5905
$ws.Range("N126").Value = -10764.0002
$ws.Range("H136").Value = 417867.88
$ws.Range("I136").Value = 1386.4
$ws.Range("K136").Value = 4159.200000000001
$ws.Range("M136").Value = -1609.200000000001
